$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 11:52"

# Espana (row 5)
$ws.Range("B5").Value = 161852
$ws.Range("C5").Value = 3579
$ws.Range("D5").Value = 59109
$ws.Range("E5").Value = 86390
$ws.Range("G5").Value = 272
$ws.Range("H5").Value = 16353

# Suiza (row 14)
$ws.Range("B14").Value = 24657
$ws.Range("C14").Value = 106
$ws.Range("E14").Value = 12554
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 1003

# Noruega (row 28)
$ws.Range("B28").Value = 6360
$ws.Range("C28").Value = 46
$ws.Range("E28").Value = 6215
$ws.Range("F28").Value = 67

# Japon (row 31)
$ws.Range("D31").Value = 762
$ws.Range("E31").Value = 5144

# Marruecos overtakes Croacia in total-cases ranking, so the sorted table
# swaps their rows (row 60 <-> row 61). Croacia's own figures are unchanged.
# Row 60 becomes Marruecos with its updated stats:
$ws.Range("A60").Value = "Marruecos"
$ws.Range("B60").Value = 1527
$ws.Range("C60").Value = 79
$ws.Range("D60").Value = 141
$ws.Range("E60").Value = 1276
$ws.Range("F60").Value = 1
$ws.Range("G60").Value = 3
$ws.Range("H60").Value = 110

# Row 61 becomes Croacia, carrying the same figures it always had:
$ws.Range("A61").Value = "Croacia"
$ws.Range("B61").Value = 1495
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 231
$ws.Range("E61").Value = 1243
$ws.Range("F61").Value = 34
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 21

# Kazajistan (row 75)
$ws.Range("D75").Value = 79
$ws.Range("E75").Value = 770

# Albania (row 96)
$ws.Range("B96").Value = 433
$ws.Range("C96").Value = 17
$ws.Range("D96").Value = 197
$ws.Range("E96").Value = 213

# Brunei (row 125)
$ws.Range("D125").Value = 104
$ws.Range("E125").Value = 31
$ws.Range("F125").Value = 2
